# Refresh the cryptocurrency Price / Volume(1h) table with the latest scrape
# (scheduled GitHub Actions update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text numbers (e.g. "216.70"). Assigning a bare
# numeric-looking string lets Excel auto-convert the cell to a Number,
# which would silently drop meaningful trailing zeros (216.70 -> 216.7).
# Prefix with a leading apostrophe - exactly like typing it into Excel -
# to force the cell to stay Text.
function Set-TextValue($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
}

$ws.Range("D2").Value = '26.496.65'
$ws.Range("E2").Value = '  +2.60%  '
$ws.Range("D3").Value = '1.678.51'
$ws.Range("E3").Value = '  +3.62%  '
$ws.Range("E4").Value = '  +0.04%  '
Set-TextValue "D5" '216.70'
$ws.Range("E5").Value = '  +3.70%  '
Set-TextValue "D6" '0.5321'
$ws.Range("E6").Value = '  +2.34%  '
$ws.Range("E7").Value = '  +0.02%  '
Set-TextValue "D8" '0.2676'
$ws.Range("E8").Value = '  +4.37%  '
Set-TextValue "D9" '0.06396'
$ws.Range("E9").Value = '  +1.68%  '
Set-TextValue "D10" '21.63'
$ws.Range("E10").Value = '  +6.10%  '
Set-TextValue "D11" '0.07806'
$ws.Range("E11").Value = '  +3.82%  '
$ws.Range("D12").Value = '1.683.14'
$ws.Range("E12").Value = '  +3.84%  '
Set-TextValue "D13" '4.499'
$ws.Range("E13").Value = '  +3.05%  '
Set-TextValue "D14" '0.5568'
$ws.Range("E14").Value = '  +2.29%  '
$ws.Range("D15").Value = '0.0₅8366'
$ws.Range("E15").Value = '  +5.44%  '
Set-TextValue "D16" '65.75'
$ws.Range("E16").Value = '  +2.66%  '
$ws.Range("D17").Value = '26.525.56'
$ws.Range("E17").Value = '  +2.69%  '
Set-TextValue "D18" '1.001'
$ws.Range("E18").Value = '  +0.06%  '
Set-TextValue "D19" '4.771'
$ws.Range("E19").Value = '  +2.81%  '
Set-TextValue "D20" '194.84'
$ws.Range("E20").Value = '  +6.22%  '
$ws.Range("E21").Value = '  +3.28%  '
Set-TextValue "D22" '6.329'
$ws.Range("E22").Value = '  +4.38%  '
$ws.Range("E23").Value = '  +0.08%  '
Set-TextValue "D24" '143.80'
$ws.Range("E24").Value = '  -0.48%  '
Set-TextValue "D25" '0.1281'
$ws.Range("E25").Value = '  +6.46%  '
Set-TextValue "D26" '7.436'
$ws.Range("E26").Value = '  +1.13%  '
Set-TextValue "D27" '16.30'
$ws.Range("E27").Value = '  +5.06%  '
Set-TextValue "D28" '1.429'
$ws.Range("E28").Value = '  +5.50%  '
Set-TextValue "D29" '0.06138'
$ws.Range("E29").Value = '  +4.61%  '
$ws.Range("E30").Value = '  +2.88%  '
$ws.Range("E31").Value = '  +7.34%  '
Set-TextValue "D32" '3.453'
$ws.Range("E32").Value = '  +3.35%  '
Set-TextValue "D33" '1.691'
$ws.Range("E33").Value = '  +5.16%  '
Set-TextValue "D34" '1.007'
$ws.Range("E34").Value = '  +3.77%  '
Set-TextValue "D35" '2.425'
$ws.Range("E35").Value = '  +1.89%  '
Set-TextValue "D36" '2.782'
$ws.Range("E36").Value = '  +2.19%  '
Set-TextValue "D37" '0.5747'
$ws.Range("E37").Value = '  -0.20%  '
Set-TextValue "D38" '0.01639'
$ws.Range("E38").Value = '  +3.46%  '
Set-TextValue "D39" '6.047'
$ws.Range("E39").Value = '  +6.91%  '
$ws.Range("D40").Value = '1.074.91'
$ws.Range("E40").Value = '  +5.38%  '
Set-TextValue "D41" '0.8605'
$ws.Range("E41").Value = '  +2.29%  '
Set-TextValue "D42" '0.9999'
$ws.Range("E42").Value = '  -0.18%  '
Set-TextValue "D43" '100.03'
$ws.Range("E43").Value = '  +0.73%  '
$ws.Range("D44").Value = '1.826.26'
$ws.Range("E44").Value = '  +3.46%  '
Set-TextValue "D45" '57.08'
$ws.Range("E45").Value = '  +5.03%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D46" '8.154'
$ws.Range("E46").Value = '  +2.79%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.0₈104'
$ws.Range("E47").Value = '  -4.75%  '
$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue "D48" '1.003'
$ws.Range("E48").Value = '  +0.11%  '
Set-TextValue "D50" '1.474'
$ws.Range("E50").Value = '  +7.14%  '
Set-TextValue "D51" '6.038'
$ws.Range("E51").Value = '  +4.01%  '
